$wb = $excel.ActiveWorkbook

# Rename the language sheet from "language_English" to "language_ENG"
$langSheet = $wb.Worksheets.Item("language_English")
$langSheet.Name = "language_ENG"

# Update the "default language name" value on the settings sheet from "English" to "ENG"
$settings = $wb.Worksheets.Item("settings")
$settings.Activate()
$settings.Range("A2").Select()
$settings.Range("A2").Value = "ENG"

# Make "tags" the active sheet (activeTab = 0)
$tags = $wb.Worksheets.Item("tags")
$tags.Activate()
